# Update task list in project
$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Task list" sheet ---
$ws = $wb.Worksheets.Item("Task list")
$ws.Name = "Task list - Core library"

# --- 2. Add new row 4 with a new task (fill left-to-right so the shared
#        string table is populated in the same order as the target file) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "normal"

$descC = "TextDocument should handle duplicated TextRanges"
$ws.Range("C4").Value = $descC
$ws.Range("C4").Characters(1, 12).Font.Bold = $true
$posRanges = $descC.IndexOf("TextRanges") + 1
$ws.Range("C4").Characters($posRanges, 10).Font.Bold = $true

$descD = "When we register or create a new TextRange then TextDocument should check its existence to avoid duplications"
$ws.Range("D4").Value = $descD
$posRange = $descD.IndexOf("TextRange ") + 1
$ws.Range("D4").Characters($posRange, 9).Font.Bold = $true
$posDoc = $descD.IndexOf("TextDocument", $posRange) + 1
$ws.Range("D4").Characters($posDoc, 12).Font.Bold = $true

# --- 3. Update status cells to "done" (row 3 existing task + new row 4) ---
$ws.Range("E3").Value = "done"
$ws.Range("E4").Value = "done"

$ws.Range("F4").Value = 42041
$ws.Range("G4").Value = 42041

$ws.Rows.Item(4).RowHeight = 45

# --- 4. Update selection on the "Task list" sheet ---
[void]$ws.Range("A3:G4").Select()

# --- 5. Clear the stray "Tasks" cell on "Project Summary" sheet ---
$ws2 = $wb.Worksheets.Item("Project Summary")
$ws2.Range("A1").ClearContents()
